{"js": "// Update the date heading (first paragraph) and all the arithmetic\n// problems in the table's cells to the new values.\n\nconst newDate = \"2026-02-13 Friday\";\n\nconst newTableValues = [\n  [\"86-67=\", \"52+16=\", \"28+51=\", \"70+22=\", \"19-5=\"],\n  [\"88-20=\", \"13-6=\", \"62-5=\", \"36+7=\", \"86-16=\"],\n  [\"51-3=\", \"74-52=\", \"41+55=\", \"54-45=\", \"12-4=\"],\n  [\"45+35=\", \"8+15=\", \"72-62=\", \"78-4=\", \"93-65=\"],\n  [\"14-11=\", \"27+66=\", \"28+26=\", \"61+6=\", \"20+57=\"],\n  [\"18+43=\", \"40-29=\", \"23-17=\", \"64-10=\", \"78+14=\"],\n  [\"67+28=\", \"32-5=\", \"2+13=\", \"86-62=\", \"86-13=\"],\n  [\"76-36=\", \"93-22=\", \"17+6=\", \"62-4=\", \"9+60=\"],\n  [\"3+53=\", \"62-0=\", \"58-29=\", \"85-69=\", \"38+6=\"],\n  [\"85-72=\", \"91-85=\", \"65+27=\", \"30-2=\", \"96+1=\"],\n  [\"95-70=\", \"33+17=\", \"31+38=\", \"13+51=\", \"92+0=\"],\n  [\"93-42=\", \"87-39=\", \"50+15=\", \"69-9=\", \"24+2=\"],\n  [\"46-28=\", \"96-9=\", \"95-32=\", \"70+8=\", \"68-18=\"],\n  [\"91-87=\", \"75-25=\", \"44+34=\", \"61+3=\", \"42-39=\"],\n  [\"67-43=\", \"24+11=\", \"80-29=\", \"23-21=\", \"51+46=\"],\n  [\"39+58=\", \"52-42=\", \"38+56=\", \"85+3=\", \"43+11=\"],\n  [\"30+29=\", \"77-8=\", \"65+5=\", \"85-22=\", \"65-45=\"],\n  [\"65+33=\", \"22+39=\", \"1+66=\", \"1+21=\", \"20+3=\"],\n  [\"66+23=\", \"21+11=\", \"97-38=\", \"74-42=\", \"89-29=\"],\n  [\"1+42=\", \"59-15=\", \"26+6=\", \"51-19=\", \"90-75=\"],\n];\n\n// 1) Update the first paragraph (the date line) in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n// Replace the whole run's text with the new date, preserving formatting.\ndateParagraph.insertText(newDate, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Update every cell of the (single) table with the new arithmetic\n// problems. The table keeps the same 20 rows x 5 columns shape, so a\n// direct bulk values assignment is sufficient.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newTableValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date heading (first paragraph of the document body).\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null\n$dateRange.Text = \"2026-02-13 Friday\"\n\n# 2) Update every cell of the (single) table with the new arithmetic\n# problems. The table keeps its original 20 rows x 5 columns shape, so\n# each cell is updated in place by (row, column) position.\n$newValues = @(\n    @(\"86-67=\",\"52+16=\",\"28+51=\",\"70+22=\",\"19-5=\"),\n    @(\"88-20=\",\"13-6=\",\"62-5=\",\"36+7=\",\"86-16=\"),\n    @(\"51-3=\",\"74-52=\",\"41+55=\",\"54-45=\",\"12-4=\"),\n    @(\"45+35=\",\"8+15=\",\"72-62=\",\"78-4=\",\"93-65=\"),\n    @(\"14-11=\",\"27+66=\",\"28+26=\",\"61+6=\",\"20+57=\"),\n    @(\"18+43=\",\"40-29=\",\"23-17=\",\"64-10=\",\"78+14=\"),\n    @(\"67+28=\",\"32-5=\",\"2+13=\",\"86-62=\",\"86-13=\"),\n    @(\"76-36=\",\"93-22=\",\"17+6=\",\"62-4=\",\"9+60=\"),\n    @(\"3+53=\",\"62-0=\",\"58-29=\",\"85-69=\",\"38+6=\"),\n    @(\"85-72=\",\"91-85=\",\"65+27=\",\"30-2=\",\"96+1=\"),\n    @(\"95-70=\",\"33+17=\",\"31+38=\",\"13+51=\",\"92+0=\"),\n    @(\"93-42=\",\"87-39=\",\"50+15=\",\"69-9=\",\"24+2=\"),\n    @(\"46-28=\",\"96-9=\",\"95-32=\",\"70+8=\",\"68-18=\"),\n    @(\"91-87=\",\"75-25=\",\"44+34=\",\"61+3=\",\"42-39=\"),\n    @(\"67-43=\",\"24+11=\",\"80-29=\",\"23-21=\",\"51+46=\"),\n    @(\"39+58=\",\"52-42=\",\"38+56=\",\"85+3=\",\"43+11=\"),\n    @(\"30+29=\",\"77-8=\",\"65+5=\",\"85-22=\",\"65-45=\"),\n    @(\"65+33=\",\"22+39=\",\"1+66=\",\"1+21=\",\"20+3=\"),\n    @(\"66+23=\",\"21+11=\",\"97-38=\",\"74-42=\",\"89-29=\"),\n    @(\"1+42=\",\"59-15=\",\"26+6=\",\"51-19=\",\"90-75=\")\n)\n\n$tbl = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Length; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Length; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
